# Adds a new "Game" Heading2 section with one numbered to-do item,
# appended after the last paragraph of the document ("Start action
# system, give nonstatic objects the ability to do different stuff on
# call ").

$d = $word.ActiveDocument

# Locate the very end of the document's main story and get a
# zero-length (collapsed) range there. Using a genuinely empty range
# makes InsertXML append new content instead of replacing existing
# paragraph content.
$endPos = $d.Content.End
$insertionPoint = $d.Range($endPos, $endPos)

$wordMlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParagraphsXml =
  "<w:p $wordMlNs>" +
    "<w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr>" +
    "<w:r><w:t>Game</w:t></w:r>" +
  "</w:p>" +
  "<w:p $wordMlNs>" +
    "<w:pPr>" +
      "<w:pStyle w:val=`"ListParagraph`"/>" +
      "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr>" +
    "</w:pPr>" +
    "<w:r><w:t>Add globals references to keyboard keys</w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`"> and implement into player controls</w:t></w:r>" +
  "</w:p>"

$insertionPoint.InsertXML($newParagraphsXml)
